# Update "想去人数" (column F) counts on several rows across the four
# worksheets, per the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 27
$ws.Cells.Item(6, 6).Value = 685
$ws.Cells.Item(7, 6).Value = 1250
$ws.Cells.Item(9, 6).Value = 851
$ws.Cells.Item(10, 6).Value = 710
$ws.Cells.Item(11, 6).Value = 266
$ws.Cells.Item(13, 6).Value = 377
$ws.Cells.Item(15, 6).Value = 995
$ws.Cells.Item(16, 6).Value = 11188
$ws.Cells.Item(17, 6).Value = 641
$ws.Cells.Item(22, 6).Value = 281
$ws.Cells.Item(23, 6).Value = 1789
$ws.Cells.Item(27, 6).Value = 189
$ws.Cells.Item(28, 6).Value = 109
$ws.Cells.Item(29, 6).Value = 285
$ws.Cells.Item(30, 6).Value = 199
$ws.Cells.Item(32, 6).Value = 77
$ws.Cells.Item(33, 6).Value = 103
$ws.Cells.Item(37, 6).Value = 191

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 6
$ws.Cells.Item(3, 6).Value = 16
$ws.Cells.Item(7, 6).Value = 144
$ws.Cells.Item(10, 6).Value = 246
$ws.Cells.Item(11, 6).Value = 4441
$ws.Cells.Item(16, 6).Value = 320

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 827

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 827
$ws.Cells.Item(3, 6).Value = 27
$ws.Cells.Item(4, 6).Value = 6
$ws.Cells.Item(7, 6).Value = 16
$ws.Cells.Item(9, 6).Value = 685
$ws.Cells.Item(10, 6).Value = 1250
$ws.Cells.Item(13, 6).Value = 144
$ws.Cells.Item(14, 6).Value = 851
$ws.Cells.Item(15, 6).Value = 710
$ws.Cells.Item(16, 6).Value = 266
$ws.Cells.Item(18, 6).Value = 995
$ws.Cells.Item(19, 6).Value = 11188
$ws.Cells.Item(20, 6).Value = 246
$ws.Cells.Item(21, 6).Value = 641
$ws.Cells.Item(24, 6).Value = 281
$ws.Cells.Item(25, 6).Value = 1789
$ws.Cells.Item(27, 6).Value = 189
$ws.Cells.Item(33, 6).Value = 320
$ws.Cells.Item(34, 6).Value = 285
$ws.Cells.Item(36, 6).Value = 199
$ws.Cells.Item(38, 6).Value = 77
$ws.Cells.Item(39, 6).Value = 103
$ws.Cells.Item(46, 6).Value = 191
